$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "sector" labels in column B (rows 2-13) with the new
# overarching "Sector Framework" labels, and clear the explicit cell
# style so these cells fall back to the default style (no "s" attribute).

$ws.Range("B2:B4").Value = "National Water Resource Management Sector Framework"
$ws.Range("B2:B4").Style = "Normal"

$ws.Range("B5:B7").Value = "National Disaster Risk Management Sector Framework"
$ws.Range("B5:B7").Style = "Normal"

$ws.Range("B8:B10").Value = "Overarching National Drought Risk Management Framework"
$ws.Range("B8:B10").Style = "Normal"

$ws.Range("B11:B13").Value = "Overarching Flood Risk Management Framework"
$ws.Range("B11:B13").Style = "Normal"

# Update the frozen-pane scroll position (topLeftCell) and the current
# selection to match the saved view state.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B11:B13").Select()
